# Auto-generated edit script: rebuild the 'Axar Patel ' sheet table with the new ownTeam/oppTeam columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear any existing content/formatting so stale cells beyond the new range don't linger
$ws.Cells.Clear()

# Header row (row 1)
$headers = @('venue', 'date', 'result', 'ownTeam', 'oppTeam', 'batsman', 'totalRuns', 'totalBalls', 'total4s', 'total6s', 'sr')
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows (rows 2..10)
$rows = @(
    @(' Dubai (DSC)', ' October 27 2020', 'Sunrisers won by 88 runs', 'Delhi Capitals', 'Sunrisers Hyderabad', 'Axar Patel ', '''1', '''4', '''0', '''0', '''25.00'),
    @(' Dubai (DSC)', ' November 05 2020', 'Mumbai won by 57 runs', 'Delhi Capitals', 'Mumbai Indians', 'Axar Patel ', '''42', '''33', '''2', '''3', '''127.27'),
    @(' Dubai (DSC)', ' November 10 2020', 'Mumbai won by 5 wickets (with 8 balls remaining)', 'Delhi Capitals', 'Mumbai Indians', 'Axar Patel ', '''9', '''9', '''1', '''0', '''100.00'),
    @(' Dubai (DSC)', ' October 14 2020', 'Capitals won by 13 runs', 'Delhi Capitals', 'Rajasthan Royals', 'Axar Patel ', '''7', '''4', '''1', '''0', '''175.00'),
    @(' Sharjah', ' October 17 2020', 'Capitals won by 5 wickets (with 1 ball remaining)', 'Delhi Capitals', 'Chennai Super Kings', 'Axar Patel ', '''21', '''5', '''0', '''3', '''420.00'),
    @(' Sharjah', ' October 09 2020', 'Capitals won by 46 runs', 'Delhi Capitals', 'Rajasthan Royals', 'Axar Patel ', '''17', '''8', '''2', '''1', '''212.50'),
    @(' Dubai (DSC)', ' September 20 2020', 'Match tied (Capitals won the one-over eliminator)', 'Delhi Capitals', 'Kings XI Punjab', 'Axar Patel ', '''6', '''9', '''0', '''0', '''66.66'),
    @(' Abu Dhabi', ' October 24 2020', 'KKR won by 59 runs', 'Delhi Capitals', 'Kolkata Knight Riders', 'Axar Patel ', '''9', '''7', '''0', '''1', '''128.57'),
    @(' Abu Dhabi', ' September 29 2020', 'Sunrisers won by 15 runs', 'Delhi Capitals', 'Sunrisers Hyderabad', 'Axar Patel ', '''5', '''6', '''0', '''0', '''83.33'),
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

Write-Output "wrote headers + $($rows.Length) data rows"
